$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FilesTab query (row 5, column B) - add LIMIT 100
$filesQuery = @'
with file_data as (
select file_name, data_category,file_type, file_size,file_access,  file_description,"sample.id" from df_sequencing_file
union
select file_name, data_category,file_type, file_size,file_access,  file_description,"sample.id" from df_methylation_array_file)

SELECT DISTINCT
    fd.file_name AS "File Name",
    fd.data_category AS "Data Category",
    COALESCE(fd.file_description, '') AS "File Description",
    fd.file_type AS "File Type",
    CASE     
        WHEN fd.file_size >= 1024 * 1024 * 1024 THEN 
            ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB' 
        WHEN fd.file_size >= 1024 * 1024 THEN 
            ROUND(fd.file_size / (1024.0 * 1024.0), 2) || ' MB' 
        WHEN fd.file_size >= 1024 THEN 
            ROUND(fd.file_size / 1024.0, 2) || ' KB' 
        ELSE 
            ROUND(fd.file_size, 2) || ' Bytes' 
    END AS "File Size",
    fd.file_access AS "File Access",
    std.dbgap_accession AS "Study ID",
    prt.participant_id AS "Participant ID",
    smp.sample_id AS "Sample ID"    
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
JOIN 
    file_data fd ON smp.id = fd."sample.id"
WHERE 
    std.dbgap_accession = 'phs002504'
LIMIT 100;
'@
$ws.Range("B5").Value = $filesQuery

# SamplesTab query (row 4, column B) - add LIMIT 100
$samplesQuery = @'
SELECT DISTINCT
   smp.sample_id AS "Sample ID",
    prt.participant_id AS "Participant ID", std.dbgap_accession AS "Study ID" , std.id, smp.anatomic_site AS "Sample Anatomic Site", dgn. "participant.id",
    COALESCE(CASE WHEN smp.participant_age_at_collection = -999 THEN 'Not Reported' ELSE smp.participant_age_at_collection END, 0) AS "Age at Sample Collection (days)",
    COALESCE(smp.sample_tumor_status, '') AS "Sample Tumor Status",
    COALESCE(smp.tumor_classification, '') AS "Sample Tumor Classification",
--COALESCE(CASE WHEN dgn."participant.id" is null THEN dgn.diagnosis ELSE dgn.diagnosis_comment END,dgn.diagnosis_comment) AS "Sample Diagnosis1",
 dgn.diagnosis as "Sample Diagnosis"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON smp."id" = dgn."sample.id"
WHERE 
   std.dbgap_accession = 'phs002504'
AND smp.sample_id IS NOT NULL
ORDER BY 
   smp.sample_id ASC
LIMIT 100;
'@
$ws.Range("B4").Value = $samplesQuery

# ParticipantsTab query (row 2, column B) - add LIMIT 100
$participantsQuery = @'
with diagnosis1 as (
select dgn."participant.id", group_concat(dgn.age_at_diagnosis,';') as age, group_concat(dgn.diagnosis,';') as diag,group_concat(dgn.anatomic_site,';') as ant_site from df_diagnosis dgn where dgn."participant.id" is not null group by dgn."participant.id" ),
diagnosis2 as (select "participant.id",  group_concat(diagnosis,';') as diag from (select distinct "participant.id", diagnosis from df_diagnosis  where "participant.id" is not null )  group by "participant.id" ),
diagnosis3 as (select "participant.id",  group_concat(anatomic_site,';') as ant_site from (select distinct "participant.id", anatomic_site from df_diagnosis where "participant.id" is not null ) group by "participant.id" ),
treatment1 as (select trt."participant.id",trt.treatment_type from  df_treatment trt  where trt.treatment_type is not null)
SELECT DISTINCT
    prt.participant_id AS "Participant ID",
    std.dbgap_accession AS "Study ID",
    COALESCE(prt.sex_at_birth, '') AS "Sex",
    COALESCE(prt.race, '') AS "Race",
	dgn2.diag AS "Diagnosis",
	dgn3.ant_site AS "Diagnosis Anatomic Site",
	   COALESCE(CASE WHEN dgn1.age = '-999' THEN 'Not Reported' ELSE dgn1.age END, "") AS "Age at Diagnosis (days)",
	trt1.treatment_type AS "Treatment Type",
	srv.last_known_survival_status AS "Last Known Survival Status"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    diagnosis1 dgn1 ON prt.id = dgn1."participant.id" 
LEFT JOIN 
    diagnosis2 dgn2 ON prt.id = dgn2."participant.id"
LEFT JOIN 
    diagnosis3 dgn3 ON prt.id = dgn3."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    treatment1 trt1 ON prt.id = trt1."participant.id"
WHERE 
    std.dbgap_accession = 'phs002504'
ORDER BY 
    prt.participant_id ASC 
LIMIT 100;
'@
$ws.Range("B2").Value = $participantsQuery

# Update the active selection to match the saved view (cell B4)
$ws.Range("B4").Select()

